# CoownersTable rewrite method searchByUserField - now its universal method
# to search by any user Field (refactor of the test workbook row used by
# the test). Update the "tttttt" test-fixture values in row 11 for the
# first_Name / login / middle_name / password columns to the new "ttttt"
# probe value, then scroll/select the sheet the way it was left after the
# edit (view scrolled to column J, S2:S11 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edit: row 11, columns H (first_Name), L (login), M (middle_name),
#     N (password) go from "tttttt" to the new probe value "ttttt".
#     (I11/last_modified, J11/last_Name and K11/locked_till stay untouched.)
$ws.Range("H11").Value = "ttttt"
$ws.Range("L11").Value = "ttttt"
$ws.Range("M11").Value = "ttttt"
$ws.Range("N11").Value = "ttttt"

# --- view state: scroll the window so column J is the left-most visible
#     column, then select S2:S11 with S2 as the active cell.
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 10   # column J
$win.ScrollRow = 1

[void]$ws.Range("S2:S11").Select()
